# Updates the cryptos list with newly scraped prices / volume percentages.
# Values in column D are stored as plain text (not numbers) in the workbook,
# so an apostrophe prefix is used to force Excel to treat them as text
# instead of auto-converting to a numeric value; the style is then reset
# back to Normal so no extra "quote prefix" formatting is retained.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'43.666.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.329.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.50%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'268.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6 - Solana
$ws.Range("E6").Value = "  +6.65%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +2.08%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'44.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.51%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.43%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'7.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.16%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.31%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "'2.669.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.05%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'15.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.85%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "'0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.50%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'2.317.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.43%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'43.621.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +1.84%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'6.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.42%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "'70.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.50%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "'239.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.67%  "

# Row 23 - ImmutableX
$ws.Range("E23").Value = "  -4.21%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'9.62"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.42%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.08%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  -7.37%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  +3.39%  "

# Row 28 - WEMIXToken
$ws.Range("E28").Value = "  -4.55%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.33%  "

# Rows 30/31 - InjectiveProtocol and EthereumClassic swap places (rank unchanged)
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.19%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'38.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "

# Row 32 - Monero
$ws.Range("D32").Value = "'172.45"
$ws.Range("D32").Style = "Normal"

# Row 33 - Hedera
$ws.Range("E33").Value = "  +0.04%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "'5.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.39%  "

# Row 35 - Stellar
$ws.Range("E35").Value = "  +1.49%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  +0.15%  "

# Row 37 - RenderToken
$ws.Range("D37").Value = "'4.45"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.47%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  -2.49%  "

# Row 39 - NEARProtocol
$ws.Range("D39").Value = "'3.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.86%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").Value = "'2.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.50%  "

# Row 41 - Algorand
$ws.Range("E41").Value = "  +15.67%  "

# Row 42 - ARBITRUM
$ws.Range("D42").Value = "'1.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.98%  "

# Row 43 - Celestia
$ws.Range("D43").Value = "'12.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "

# Row 44 - THORChain
$ws.Range("D44").Value = "'5.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "

# Row 45 - FraxShare
$ws.Range("D45").Value = "'9.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.63%  "

# Row 46 - MultiversX
$ws.Range("D46").Value = "'60.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.95%  "

# Row 47 - Cronos
$ws.Range("E47").Value = "  +2.90%  "

# Row 48 - Aave
$ws.Range("D48").Value = "'100.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.61%  "

# Row 49 - TrustWalletToken
$ws.Range("D49").Value = "'1.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "'2.551.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.16%  "

# Row 51 - WOONetwork
$ws.Range("E51").Value = "  -3.04%  "
